# Applies the data edits from the commit: a few cell values on Sheet1
# changed (row 3: E3, G3, H3) and the active selection moved from A6 to E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 5
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 13

$ws.Range("E3").Select()
